$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-location")

# Data rows to append below the current last row (row 109 -> new rows 110..119)
$rows = @(
    @(10110, 10110, 5, "Postal Code", "BNMR", "eng"),
    @(10111, 10111, 5, "Postal Code", "BNMR", "eng"),
    @(10113, 10113, 5, "Postal Code", "BNMR", "eng"),
    @(10114, 10114, 5, "Postal Code", "BNMR", "eng"),
    @(10111, 10111, 5, "code postal", "BNMR", "fra"),
    @(10110, 10110, 5, "code postal", "BNMR", "fra"),
    @(10113, 10113, 5, "code postal", "BNMR", "fra"),
    @(10114, 10114, 5, "code postal", "BNMR", "fra"),
    @(10111, 10111, 5, "الرمز البريدي", "BNMR", "ara"),
    @(10110, 10110, 5, "الرمز البريدي", "BNMR", "ara")
)

$startRow = 110
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $true
    $ws.Cells.Item($r, 8).Value = "superadmin"
    $ws.Cells.Item($r, 9).Value = "now()"
}
